$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").Value = "TrackingError1"
$ws.Range("B60").Value = "Test Tracking Error with scale=252"
$ws.Range("C60").Value = "Tracking_Error_test1"

$ws.Range("A61").Value = "TrackingError2"
$ws.Range("B61").Value = "Test Tracking Error with scale=1"
$ws.Range("C61").Value = "Tracking_Error_test2"

$ws.Range("G64").Select() | Out-Null
